# semana 40 de 2025
# Adds the epidemiological-week 39 and 40 columns (AP, AQ) to the weekly
# IRA-hospitalario report: header labels in row 1, plus the per-facility
# counts that were reported for those two weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns AP (week 39) and AQ (week 40), matching the
# formatting already used by the rest of the header row (row 1): bold,
# centered (same as style index 1 applied to A1:AO1).
$ws.Range("AP1").Value = "39"
$ws.Range("AQ1").Value = "40"
$ws.Range("AP1:AQ1").Font.Bold = $true
$ws.Range("AP1:AQ1").HorizontalAlignment = -4108

# Per-row counts for weeks 39 (AP) and 40 (AQ). Rows not listed here (e.g.
# 4, 18-21, 27, 32-34, 39, 40) received no new data for these weeks.
$weekData = @{
    2  = @{ "AP" = 0;  "AQ" = 0 }
    3  = @{ "AP" = 0;  "AQ" = 0 }
    5  = @{ "AP" = 0;  "AQ" = 0 }
    6  = @{ "AP" = 28; "AQ" = 31 }
    7  = @{ "AP" = 4;  "AQ" = 2 }
    8  = @{ "AP" = 26; "AQ" = 22 }
    9  = @{ "AP" = 0;  "AQ" = 0 }
    10 = @{ "AP" = 0;  "AQ" = 0 }
    11 = @{ "AP" = 0 }
    12 = @{ "AP" = 0;  "AQ" = 0 }
    13 = @{ "AP" = 0 }
    14 = @{ "AP" = 0;  "AQ" = 0 }
    15 = @{ "AQ" = 0 }
    16 = @{ "AP" = 0;  "AQ" = 0 }
    17 = @{ "AP" = 0;  "AQ" = 0 }
    22 = @{ "AP" = 0;  "AQ" = 0 }
    23 = @{ "AP" = 0;  "AQ" = 0 }
    24 = @{ "AQ" = 0 }
    25 = @{ "AP" = 1;  "AQ" = 1 }
    26 = @{ "AP" = 0 }
    28 = @{ "AP" = 3;  "AQ" = 1 }
    29 = @{ "AP" = 3;  "AQ" = 1 }
    30 = @{ "AP" = 1;  "AQ" = 0 }
    31 = @{ "AP" = 0;  "AQ" = 0 }
    35 = @{ "AP" = 4;  "AQ" = 8 }
    36 = @{ "AP" = 1;  "AQ" = 0 }
    37 = @{ "AP" = 0;  "AQ" = 0 }
    38 = @{ "AP" = 0;  "AQ" = 0 }
    41 = @{ "AP" = 0;  "AQ" = 0 }
    42 = @{ "AP" = 0;  "AQ" = 0 }
    43 = @{ "AP" = 0;  "AQ" = 0 }
    44 = @{ "AP" = 0 }
    45 = @{ "AP" = 0;  "AQ" = 0 }
    46 = @{ "AP" = 0;  "AQ" = 0 }
    47 = @{ "AP" = 0;  "AQ" = 0 }
    48 = @{ "AP" = 0;  "AQ" = 0 }
    49 = @{ "AP" = 0;  "AQ" = 0 }
    50 = @{ "AP" = 0;  "AQ" = 0 }
    51 = @{ "AP" = 0;  "AQ" = 0 }
    52 = @{ "AP" = 0;  "AQ" = 0 }
    53 = @{ "AP" = 0;  "AQ" = 0 }
    54 = @{ "AP" = 0;  "AQ" = 0 }
    55 = @{ "AP" = 0;  "AQ" = 0 }
    56 = @{ "AP" = 0;  "AQ" = 0 }
    57 = @{ "AP" = 0;  "AQ" = 0 }
    58 = @{ "AP" = 0;  "AQ" = 0 }
}

foreach ($rowNum in $weekData.Keys) {
    $cols = $weekData[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
